# SOEN 6011 Skill Selection Sheet - add skill entry for Dhruv Goyani (row 9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Skill" column value for Dhruv Goyani's row from
# "Algorithms and Complexity" to "Algorithms and Complexity Analysis"
$ws.Range("C9").Value = "Algorithms and Complexity Analysis"

# Leave the active selection on C9, matching where the edit was made
$ws.Range("C9").Select() | Out-Null
